$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2027491408934708
$ws.Range("C2").Value = 0.5532646048109966
$ws.Range("J2").Value = 0.01030927835051546
$ws.Range("P2").Value = 0.1512027491408935
$ws.Range("S2").Value = 0.08247422680412371
$ws.Range("B3").Value = 0.0119047619047619
$ws.Range("C3").Value = 0.04166666666666666
$ws.Range("J3").Value = 0.02380952380952381
$ws.Range("P3").Value = 0.7202380952380952
$ws.Range("S3").Value = 0.2023809523809524
$ws.Range("J4").Value = 0.07142857142857142
$ws.Range("P4").Value = 0.5952380952380952
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.05376344086021505
$ws.Range("D6").Value = 0.01075268817204301
$ws.Range("F6").Value = 0.07526881720430108
$ws.Range("J6").Value = 0.2329749103942652
$ws.Range("O6").Value = 0.02867383512544803
$ws.Range("Q6").Value = 0.1505376344086022
$ws.Range("R6").Value = 0.05734767025089606
$ws.Range("S6").Value = 0.3906810035842294
$ws.Range("B7").Value = 0.1063829787234043
$ws.Range("D7").Value = 0.02127659574468085
$ws.Range("F7").Value = 0.07234042553191489
$ws.Range("J7").Value = 0.1531914893617021
$ws.Range("O7").Value = 0.02553191489361702
$ws.Range("Q7").Value = 0.1191489361702128
$ws.Range("R7").Value = 0.05106382978723404
$ws.Range("S7").Value = 0.451063829787234
$ws.Range("B8").Value = 0.09876543209876543
$ws.Range("D8").Value = 0.01234567901234568
$ws.Range("E8").Value = 0.00205761316872428
$ws.Range("F8").Value = 0.06995884773662552
$ws.Range("J8").Value = 0.08024691358024691
$ws.Range("O8").Value = 0.02674897119341564
$ws.Range("Q8").Value = 0.1790123456790123
$ws.Range("R8").Value = 0.08847736625514403
$ws.Range("S8").Value = 0.4423868312757201
$ws.Range("B9").Value = 0.07627118644067797
$ws.Range("D9").Value = 0.0211864406779661
$ws.Range("F9").Value = 0.07627118644067797
$ws.Range("J9").Value = 0.1016949152542373
$ws.Range("O9").Value = 0.05084745762711865
$ws.Range("Q9").Value = 0.1525423728813559
$ws.Range("R9").Value = 0.1016949152542373
$ws.Range("S9").Value = 0.4194915254237288
$ws.Range("B10").Value = 0.09043736100815419
$ws.Range("D10").Value = 0.01704966641957005
$ws.Range("F10").Value = 0.07635285396590066
$ws.Range("J10").Value = 0.1111934766493699
$ws.Range("O10").Value = 0.01779095626389918
$ws.Range("Q10").Value = 0.2127501853224611
$ws.Range("R10").Value = 0.08154188287620459
$ws.Range("S10").Value = 0.3928836174944403
$ws.Range("F11").Value = 0.002770083102493075
$ws.Range("G11").Value = 0.1662049861495845
$ws.Range("J11").Value = 0.08310249307479224
$ws.Range("K11").Value = 0.2022160664819945
$ws.Range("L11").Value = 0.5318559556786704
$ws.Range("S11").Value = 0.01385041551246537
$ws.Range("G12").Value = 0.735
$ws.Range("J12").Value = 0.19
$ws.Range("K12").Value = 0.02
$ws.Range("L12").Value = 0.035
$ws.Range("S12").Value = 0.02
$ws.Range("G13").Value = 0.66
$ws.Range("J13").Value = 0.3
$ws.Range("S13").Value = 0.04
$ws.Range("F15").Value = 0.02298850574712644
$ws.Range("H15").Value = 0.1417624521072797
$ws.Range("I15").Value = 0.06896551724137931
$ws.Range("J15").Value = 0.3333333333333333
$ws.Range("K15").Value = 0.05363984674329502
$ws.Range("M15").Value = 0.01532567049808429
$ws.Range("O15").Value = 0.06513409961685823
$ws.Range("S15").Value = 0.2988505747126437
$ws.Range("F16").Value = 0.02173913043478261
$ws.Range("H16").Value = 0.1739130434782609
$ws.Range("I16").Value = 0.05978260869565218
$ws.Range("J16").Value = 0.4076086956521739
$ws.Range("K16").Value = 0.1032608695652174
$ws.Range("M16").Value = 0.0108695652173913
$ws.Range("N16").Value = 0.005434782608695652
$ws.Range("O16").Value = 0.05978260869565218
$ws.Range("S16").Value = 0.1576086956521739
$ws.Range("F17").Value = 0.03125
$ws.Range("H17").Value = 0.1916666666666667
$ws.Range("I17").Value = 0.09791666666666667
$ws.Range("J17").Value = 0.3854166666666667
$ws.Range("K17").Value = 0.10625
$ws.Range("M17").Value = 0.0125
$ws.Range("N17").Value = 0.002083333333333333
$ws.Range("O17").Value = 0.07916666666666666
$ws.Range("S17").Value = 0.09375
$ws.Range("F18").Value = 0.004854368932038835
$ws.Range("H18").Value = 0.2184466019417476
$ws.Range("I18").Value = 0.05825242718446602
$ws.Range("J18").Value = 0.4368932038834951
$ws.Range("K18").Value = 0.1019417475728155
$ws.Range("M18").Value = 0.01941747572815534
$ws.Range("O18").Value = 0.06310679611650485
$ws.Range("S18").Value = 0.0970873786407767
$ws.Range("F19").Value = 0.01558265582655827
$ws.Range("H19").Value = 0.1897018970189702
$ws.Range("I19").Value = 0.1023035230352304
$ws.Range("J19").Value = 0.3597560975609756
$ws.Range("K19").Value = 0.1172086720867209
$ws.Range("M19").Value = 0.02439024390243903
$ws.Range("O19").Value = 0.06368563685636856
$ws.Range("S19").Value = 0.1253387533875339
